$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.188019275665283
$ws.Range("B1").Value = 3.908146619796753
$ws.Range("C1").Value = 3.84501576423645
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 3.506923198699951
